# Add season-record columns (Wins / Losses / Ties) to the right of the
# existing team-stat table.
#
# Source data did not include the team's season record before, only the
# player statistics. This adds three new columns - Wins, Losses, Ties -
# and fills them in for every player row with the team's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 48

# --- Header row -------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the look of the existing header cells (bold, bordered, centered)
# by copying the formatting from an existing header cell rather than
# re-creating a brand-new style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
# Every player on the roster shares the same team season record.
$wins = 73
$losses = 89
$ties = 1

$ws.Range("AD2:AD$lastRow").Value = $wins
$ws.Range("AE2:AE$lastRow").Value = $losses
$ws.Range("AF2:AF$lastRow").Value = $ties

$excel.CutCopyMode = $false
